# Generate Report for Handoff
#
# Updates the localization-status report:
#  - Bumps the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
#    timestamps for the files that just went through another handoff cycle
#    (rows 7, 8, 10, 11, 12, 13 across all three sheets).
#  - Sets the "Priority" column to "ht" for those same rows on the
#    language sheets (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 10, 11, 12, 13)

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-16 16:18:39"
}

# --- zh-cn sheet: "Latest Handoff Datetime" (column H) + "Priority" (column E) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("H$r").Value = "2016-08-16 16:18:34"
    $wsZhCn.Range("E$r").Value = "ht"
}

# --- de-de sheet: "Latest Handoff Datetime" (column H) + "Priority" (column E) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("H$r").Value = "2016-08-16 16:18:39"
    $wsDeDe.Range("E$r").Value = "ht"
}
